$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "x"
$ws.Range("B21").Value = "x"

$ws.Range("A1:D21").Select()
$ws.Range("F20").Select()

$ws.Columns.Item(1).ColumnWidth = 61.140625
$ws.Cells.Item(1,1).EntireRow.RowHeight = 15
